$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 2, pushing the existing row 2 data down to row 4
$ws.Rows.Item(2).Resize(2).Insert()
# Excel copies formatting from the row above (header) by default; clear it so the
# new data rows have no special style, matching the rest of the data rows.
$ws.Rows.Item(2).Resize(2).ClearFormats()

# New row 2 data: OKQqJDOt / BRAZIL - SERIE B / Operario - Mirassol
$row2 = @(
    "OKQqJDOt","15/11/2024","16:00","BRAZIL - SERIE B","Operario","Mirassol",
    2.7,2.7,3,3.6,1.8,4,1.14,5.5,1.67,2.1,3.4,1.33,1.73,2,2.38,1.53,6,11,12,
    29,29,51,5,5.5,21,101,201,6.5,13,13,34,34,51,4.33,17,34,67,126,351,
    2,10,101,4.75,21,41,67,126,351,81,81
)

# New row 3 data: U9HRzAsQ / NETHERLANDS - EERSTE DIVISIE / Telstar - Venlo
$row3 = @(
    "U9HRzAsQ","15/11/2024","16:00","NETHERLANDS - EERSTE DIVISIE","Telstar","Venlo",
    1.48,4.75,5.75,1.95,2.6,5.5,1.02,19,1.14,5.5,1.48,2.6,1.25,3.75,1.67,2.1,9.5,8.5,8.5,
    11,11,21,19,10,17,41,151,19,34,17,51,41,41,3.75,7,15,19,41,81,
    3.75,8,41,7.5,29,29,101,101,151,126,151
)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}
